# dataset_summary.xlsx -- "updates to diffusion/laplacian analysis"
#
# Row 22 (the last data row, dataset "01_18_22") is updated:
#   - dataset renamed 01_18_22 -> 01_19_22
#   - k (B22) and chi (C22) columns cleared out entirely
#   - diag (E22) changes from formula =FALSE() to plain value 1
#   - N (G22) changes from 20 to 2000
#   - maxent samples (H22) cleared out entirely
#   - Comments (I22) changes from "nonlinear system of k=4" to "??"
# The selection is left on H2, matching the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E22 and I22 previously carried special formatting (a computed-boolean
# style and a bold "highlight" style respectively). The edited row goes
# back to plain/default formatting for every cell, so pull the default
# style (taken from D22, a cell in the same row that already uses it)
# onto both before changing their values.
$ws.Range("D22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("I22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# k and maxent samples columns, and chi, are removed for this row.
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("H22").ClearContents()

# Updated values for the remaining cells.
$ws.Range("A22").Value = "01_19_22"
$ws.Range("E22").Value = 1
$ws.Range("G22").Value = 2000
$ws.Range("I22").Value = "??"

# Cosmetic: the sheet's row-1048576 height was touched (a LibreOffice
# artifact) and the active selection ended on H2.
$ws.Rows.Item(1048576).RowHeight = 12.8
$ws.Range("H2").Select()
